$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 39, pushing the existing rows 39-67 down to 41-69.
$ws.Rows("39:40").Insert()

# New row 39: Black Amber / Primera, with updated volumes/prices.
$ws.Cells.Item(39, 1).Value = 5
$ws.Cells.Item(39, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value = "Maule"
$ws.Cells.Item(39, 4).Value = 44571
$ws.Cells.Item(39, 5).Value = 7
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103002
$ws.Cells.Item(39, 10).Value = "Ciruela"
$ws.Cells.Item(39, 11).Value = "Black Amber"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 230
$ws.Cells.Item(39, 14).Value = 10000
$ws.Cells.Item(39, 15).Value = 10000
$ws.Cells.Item(39, 16).Value = 10000
$ws.Cells.Item(39, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(39, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(39, 19).Value = 556
$ws.Cells.Item(39, 20).Value = 18

# New row 40: Lemon / Primera, with updated volumes/prices.
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(40, 3).Value = "Maule"
$ws.Cells.Item(40, 4).Value = 44571
$ws.Cells.Item(40, 5).Value = 7
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100103
$ws.Cells.Item(40, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value = 100103002
$ws.Cells.Item(40, 10).Value = "Ciruela"
$ws.Cells.Item(40, 11).Value = "Lemon"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 120
$ws.Cells.Item(40, 14).Value = 11000
$ws.Cells.Item(40, 15).Value = 11000
$ws.Cells.Item(40, 16).Value = 11000
$ws.Cells.Item(40, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(40, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(40, 19).Value = 611
$ws.Cells.Item(40, 20).Value = 18
